$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '29.148.13'
$ws.Cells.Item(2, 5).Value = '  -0.97%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.860.41'
$ws.Cells.Item(3, 5).Value = '  -0.68%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.08%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '0.7062'
$ws.Cells.Item(5, 5).Value = '  -0.83%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '240.90'
$ws.Cells.Item(6, 5).Value = '  -0.39%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.001'
$ws.Cells.Item(7, 5).Value = '  +0.01%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3081'
$ws.Cells.Item(8, 5).Value = '  -0.86%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.07621'
$ws.Cells.Item(9, 5).Value = '  -2.95%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '24.60'
$ws.Cells.Item(10, 5).Value = '  -1.83%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.08318'
$ws.Cells.Item(11, 5).Value = '  +0.92%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '1.863.81'
$ws.Cells.Item(12, 5).Value = '  -0.82%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  -1.73%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.7083'
$ws.Cells.Item(14, 5).Value = '  -2.42%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '90.98'
$ws.Cells.Item(15, 5).Value = '  +0.16%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '29.213.30'
$ws.Cells.Item(16, 5).Value = '  -0.81%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  -0.32%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '242.44'
$ws.Cells.Item(18, 5).Value = '  -1.76%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.000007792'
$ws.Cells.Item(19, 5).Value = '  -0.93%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '2.115.28'
$ws.Cells.Item(20, 5).Value = '  -1.07%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '13.05'
$ws.Cells.Item(21, 5).Value = '  -1.65%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '1.000'
$ws.Cells.Item(22, 5).Value = '  +0.08%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '7.846'
$ws.Cells.Item(23, 5).Value = '  -1.47%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '1.000'
$ws.Cells.Item(24, 5).Value = '  -0.10%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.1584'
$ws.Cells.Item(25, 5).Value = '  -0.37%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '163.07'
$ws.Cells.Item(26, 5).Value = '  -0.39%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '8.923'
$ws.Cells.Item(27, 5).Value = '  -0.89%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '18.42'
$ws.Cells.Item(28, 5).Value = '  +0.67%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'Toncoin'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.324'
$ws.Cells.Item(29, 5).Value = '  -3.03%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'PancakeSwap'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.497'
$ws.Cells.Item(30, 5).Value = '  +0.24%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '4.392'
$ws.Cells.Item(31, 5).Value = '  +0.46%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.199'
$ws.Cells.Item(32, 5).Value = '  +1.94%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.05125'
$ws.Cells.Item(33, 5).Value = '  -3.40%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.7938'
$ws.Cells.Item(34, 5).Value = '  +9.51%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.903'
$ws.Cells.Item(35, 5).Value = '  -1.26%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.160'
$ws.Cells.Item(36, 5).Value = '  -3.17%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.682'
$ws.Cells.Item(37, 5).Value = '  +0.19%  '

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.01837'
$ws.Cells.Item(38, 5).Value = '  -1.48%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.698'
$ws.Cells.Item(39, 5).Value = '  -0.95%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.165.35'
$ws.Cells.Item(40, 5).Value = '  -5.93%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '6.164'
$ws.Cells.Item(41, 5).Value = '  +0.17%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.8873'
$ws.Cells.Item(42, 5).Value = '  -2.22%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'Aave'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '72.77'
$ws.Cells.Item(43, 5).Value = '  -1.76%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.000'

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '102.07'
$ws.Cells.Item(45, 5).Value = '  -1.07%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.011.34'
$ws.Cells.Item(46, 5).Value = '  -0.90%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  -2.73%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.764'
$ws.Cells.Item(48, 5).Value = '  -0.48%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +0.29%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '9.278'
$ws.Cells.Item(50, 5).Value = '  -0.07%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.9997'
$ws.Cells.Item(51, 5).Value = '  -0.12%  '
